# Applies the "added harvard case classification" update:
#  - swap header labels for average_doctor / average_doctor_old (BP1/BQ1)
#  - update recomputed average/variance/stdDev statistics for rows 4-13
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header swap: the "average_doctor" / "average_doctor_old" columns traded places
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Updated statistics values

# Row 4
$ws.Range("E4").Value = 0.423
$ws.Range("F4").Value = 0.07199999999999999
$ws.Range("G4").Value = 0.268
$ws.Range("N4").Value = 0.432
$ws.Range("O4").Value = 0.064
$ws.Range("P4").Value = 0.254
$ws.Range("Q4").Value = 0.023
$ws.Range("R4").Value = 0.016
$ws.Range("S4").Value = 0.128
$ws.Range("W4").Value = 0.283
$ws.Range("X4").Value = 0.109
$ws.Range("Y4").Value = 0.33
$ws.Range("AI4").Value = 0.288
$ws.Range("AJ4").Value = 0.08799999999999999
$ws.Range("AK4").Value = 0.297
$ws.Range("AU4").Value = 0.188
$ws.Range("AV4").Value = 0.028
$ws.Range("AW4").Value = 0.168
$ws.Range("BA4").Value = 2.012
$ws.Range("BB4").Value = 0.158
$ws.Range("BC4").Value = 0.397
$ws.Range("BG4").Value = 0.732
$ws.Range("BH4").Value = 0.141
$ws.Range("BI4").Value = 0.375
$ws.Range("BM4").Value = 0.719
$ws.Range("BN4").Value = 0.08
$ws.Range("BO4").Value = 0.282
$ws.Range("BP4").Value = 0.671
$ws.Range("BQ4").Value = 0.71

# Row 5
$ws.Range("E5").Value = 0.538
$ws.Range("F5").Value = 0.08500000000000001
$ws.Range("G5").Value = 0.292
$ws.Range("N5").Value = 0.733
$ws.Range("O5").Value = 0.081
$ws.Range("P5").Value = 0.285
$ws.Range("Q5").Value = 0.015
$ws.Range("R5").Value = 0.007
$ws.Range("S5").Value = 0.082
$ws.Range("W5").Value = 0.273
$ws.Range("X5").Value = 0.109
$ws.Range("Y5").Value = 0.33
$ws.Range("AI5").Value = 0.309
$ws.Range("AJ5").Value = 0.098
$ws.Range("AK5").Value = 0.314
$ws.Range("AU5").Value = 0.365
$ws.Range("AV5").Value = 0.093
$ws.Range("AW5").Value = 0.306
$ws.Range("BA5").Value = 1.337
$ws.Range("BB5").Value = 0.079
$ws.Range("BC5").Value = 0.28
$ws.Range("BG5").Value = 0.398
$ws.Range("BH5").Value = 0.051
$ws.Range("BI5").Value = 0.225
$ws.Range("BM5").Value = 0.549
$ws.Range("BN5").Value = 0.063
$ws.Range("BO5").Value = 0.252
$ws.Range("BP5").Value = 0.446
$ws.Range("BQ5").Value = 0.456

# Row 6
$ws.Range("E6").Value = 0.474
$ws.Range("N6").Value = 0.544
$ws.Range("Q6").Value = 0.018
$ws.Range("W6").Value = 0.278
$ws.Range("AI6").Value = 0.298
$ws.Range("AU6").Value = 0.248
$ws.Range("BA6").Value = 1.599
$ws.Range("BG6").Value = 0.516
$ws.Range("BM6").Value = 0.623
$ws.Range("BP6").Value = 0.533
$ws.Range("BQ6").Value = 0.552

# Row 7
$ws.Range("E7").Value = 0.51
$ws.Range("N7").Value = 0.643
$ws.Range("Q7").Value = 0.016
$ws.Range("W7").Value = 0.275
$ws.Range("AI7").Value = 0.305
$ws.Range("AU7").Value = 0.307
$ws.Range("BA7").Value = 1.429
$ws.Range("BG7").Value = 0.438
$ws.Range("BM7").Value = 0.576
$ws.Range("BP7").Value = 0.476
$ws.Range("BQ7").Value = 0.49

# Row 8
$ws.Range("E8").Value = 0.605
$ws.Range("F8").Value = 0.112
$ws.Range("G8").Value = 0.334
$ws.Range("N8").Value = 0.773
$ws.Range("O8").Value = 0.066
$ws.Range("P8").Value = 0.256
$ws.Range("Q8").Value = 0.017
$ws.Range("S8").Value = 0.107
$ws.Range("W8").Value = 0.301
$ws.Range("Y8").Value = 0.347
$ws.Range("AI8").Value = 0.33
$ws.Range("AJ8").Value = 0.129
$ws.Range("AK8").Value = 0.359
$ws.Range("AU8").Value = 0.308
$ws.Range("AV8").Value = 0.08500000000000001
$ws.Range("AW8").Value = 0.291
$ws.Range("BA8").Value = 1.75
$ws.Range("BB8").Value = 0.124
$ws.Range("BC8").Value = 0.352
$ws.Range("BG8").Value = 0.5679999999999999
$ws.Range("BH8").Value = 0.106
$ws.Range("BI8").Value = 0.326
$ws.Range("BM8").Value = 0.6909999999999999
$ws.Range("BN8").Value = 0.066
$ws.Range("BO8").Value = 0.258
$ws.Range("BP8").Value = 0.583
$ws.Range("BQ8").Value = 0.605

# Row 9
$ws.Range("E9").Value = 0.553
$ws.Range("F9").Value = 0.247
$ws.Range("G9").Value = 0.497
$ws.Range("N9").Value = 0.681
$ws.Range("O9").Value = 0.217
$ws.Range("P9").Value = 0.466
$ws.Range("W9").Value = 0.202
$ws.Range("X9").Value = 0.161
$ws.Range("Y9").Value = 0.402
$ws.Range("AI9").Value = 0.255
$ws.Range("AJ9").Value = 0.19
$ws.Range("AK9").Value = 0.436
$ws.Range("BA9").Value = 1.712
$ws.Range("BB9").Value = 0.248
$ws.Range("BC9").Value = 0.498
$ws.Range("BG9").Value = 0.606
$ws.Range("BH9").Value = 0.239
$ws.Range("BI9").Value = 0.489
$ws.Range("BM9").Value = 0.649
$ws.Range("BN9").Value = 0.228
$ws.Range("BO9").Value = 0.477
$ws.Range("BP9").Value = 0.571
$ws.Range("BQ9").Value = 0.588

# Row 10
$ws.Range("E10").Value = 0.681
$ws.Range("F10").Value = 0.217
$ws.Range("G10").Value = 0.466
$ws.Range("N10").Value = 0.872
$ws.Range("O10").Value = 0.111
$ws.Range("P10").Value = 0.334
$ws.Range("W10").Value = 0.372
$ws.Range("X10").Value = 0.234
$ws.Range("Y10").Value = 0.483
$ws.Range("AI10").Value = 0.362
$ws.Range("AJ10").Value = 0.231
$ws.Range("AK10").Value = 0.48
$ws.Range("AU10").Value = 0.298
$ws.Range("AV10").Value = 0.209
$ws.Range("AW10").Value = 0.457
$ws.Range("BA10").Value = 2.085
$ws.Range("BB10").Value = 0.243
$ws.Range("BC10").Value = 0.493
$ws.Range("BG10").Value = 0.66
$ws.Range("BH10").Value = 0.225
$ws.Range("BI10").Value = 0.474
$ws.Range("BM10").Value = 0.84
$ws.Range("BN10").Value = 0.134
$ws.Range("BO10").Value = 0.366
$ws.Range("BP10").Value = 0.695
$ws.Range("BQ10").Value = 0.726

# Row 11
$ws.Range("E11").Value = 0.713
$ws.Range("F11").Value = 0.205
$ws.Range("G11").Value = 0.452
$ws.Range("N11").Value = 0.894
$ws.Range("O11").Value = 0.095
$ws.Range("P11").Value = 0.308
$ws.Range("W11").Value = 0.372
$ws.Range("X11").Value = 0.234
$ws.Range("Y11").Value = 0.483
$ws.Range("AI11").Value = 0.394
$ws.Range("AJ11").Value = 0.239
$ws.Range("AK11").Value = 0.489
$ws.Range("AU11").Value = 0.436
$ws.Range("AV11").Value = 0.246
$ws.Range("AW11").Value = 0.496
$ws.Range("BA11").Value = 2.085
$ws.Range("BB11").Value = 0.243
$ws.Range("BC11").Value = 0.493
$ws.Range("BG11").Value = 0.66
$ws.Range("BH11").Value = 0.225
$ws.Range("BI11").Value = 0.474
$ws.Range("BM11").Value = 0.84
$ws.Range("BN11").Value = 0.134
$ws.Range("BO11").Value = 0.366
$ws.Range("BP11").Value = 0.695
$ws.Range("BQ11").Value = 0.728

# Row 12
$ws.Range("E12").Value = 1.403
$ws.Range("F12").Value = 0.748
$ws.Range("G12").Value = 0.865
$ws.Range("N12").Value = 1.465
$ws.Range("O12").Value = 1.039
$ws.Range("P12").Value = 1.02
$ws.Range("W12").Value = 1.629
$ws.Range("X12").Value = 0.576
$ws.Range("Y12").Value = 0.759
$ws.Range("AI12").Value = 1.703
$ws.Range("AJ12").Value = 1.29
$ws.Range("AK12").Value = 1.136
$ws.Range("AU12").Value = 2.767
$ws.Range("AV12").Value = 2.737
$ws.Range("AW12").Value = 1.654
$ws.Range("BB12").Value = 0.395
$ws.Range("BC12").Value = 0.629
$ws.Range("BG12").Value = 1.097
$ws.Range("BH12").Value = 0.12
$ws.Range("BI12").Value = 0.346
$ws.Range("BM12").Value = 1.291
$ws.Range("BN12").Value = 0.333
$ws.Range("BO12").Value = 0.577
$ws.Range("BP12").Value = 1.232
$ws.Range("BQ12").Value = 1.257

# Row 13
$ws.Range("E13").Value = 1.573
$ws.Range("F13").Value = 0.652
$ws.Range("G13").Value = 0.8080000000000001
$ws.Range("N13").Value = 2.068
$ws.Range("O13").Value = 0.928
$ws.Range("P13").Value = 0.964
$ws.Range("W13").Value = 1.037
$ws.Range("X13").Value = 0.193
$ws.Range("Y13").Value = 0.439
$ws.Range("AI13").Value = 1.28
$ws.Range("AJ13").Value = 0.37
$ws.Range("AK13").Value = 0.608
$ws.Range("AU13").Value = 2.285
$ws.Range("AV13").Value = 0.925
$ws.Range("AW13").Value = 0.962
$ws.Range("BA13").Value = 2.344
$ws.Range("BB13").Value = 0.296
$ws.Range("BC13").Value = 0.544
$ws.Range("BG13").Value = 0.584
$ws.Range("BH13").Value = 0.07099999999999999
$ws.Range("BI13").Value = 0.266
$ws.Range("BM13").Value = 0.892
$ws.Range("BN13").Value = 0.282
$ws.Range("BO13").Value = 0.531
$ws.Range("BP13").Value = 0.781
$ws.Range("BQ13").Value = 0.725
